# Auto-applied edit: adds rows to "Pit Stop" and "Dry Season", and adds a new
# worksheet "Sneaky Fields" with its own header + data rows, mirroring the
# other per-arena scrim-tracking sheets already in the workbook.

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($ws, $row, $values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- "Pit Stop": append one new submission row (row 14) ---
$pitStop = $wb.Worksheets.Item("Pit Stop")
$pitStop.Range("A13:N13").Copy()
$pitStop.Range("A14:N14").PasteSpecial(-4122)
Set-RowValues $pitStop 14 @("SHADE", "BULL", "EDGAR", "NITA", "BONNIE", "LOU", "Equipo 1", "Enraged 💔", "SUP|Filippo神", "SUP|Tomzy", "IC|Mebius", "IC|RamaZR", "IC|Nob", "20250723T174004.000Z")

# --- "Dry Season": append new submission rows (rows 19-29) ---
$drySeason = $wb.Worksheets.Item("Dry Season")
$drySeason.Range("A18:N18").Copy()
$drySeason.Range("A19:N29").PasteSpecial(-4122)
Set-RowValues $drySeason 19 @("BELLE", "SQUEAK", "CARL", "OLLIE", "KAZE", "BYRON", "Equipo 2", "NXT|amos", "NXT|Arthur", "NXT|Rup", "FUT|Nowy297", "FUT|MeOw", "FUT|GeRo", "20250723T174629.000Z")
Set-RowValues $drySeason 20 @("BELLE", "SQUEAK", "CARL", "OLLIE", "KAZE", "BYRON", "Equipo 2", "NXT|amos", "NXT|Arthur", "NXT|Rup", "FUT|Nowy297", "FUT|MeOw", "FUT|GeRo", "20250723T174422.000Z")
Set-RowValues $drySeason 21 @("BELLE", "SQUEAK", "CARL", "OLLIE", "KAZE", "BYRON", "Equipo 1", "NXT|amos", "NXT|Arthur", "NXT|Rup", "FUT|Nowy297", "FUT|MeOw", "FUT|GeRo", "20250723T174223.000Z")
Set-RowValues $drySeason 22 @("MR. P", "GENE", "BELLE", "BROCK", "GUS", "BYRON", "Equipo 2", "NXT|Arthur", "NXT|Rup", "NXT|amos", "FUT|GeRo", "FUT|MeOw", "FUT|Nowy297", "20250723T173623.000Z")
Set-RowValues $drySeason 23 @("MR. P", "GENE", "BELLE", "BROCK", "GUS", "BYRON", "Equipo 1", "NXT|Arthur", "NXT|Rup", "NXT|amos", "FUT|GeRo", "FUT|MeOw", "FUT|Nowy297", "20250723T173403.000Z")
Set-RowValues $drySeason 24 @("MR. P", "GENE", "BELLE", "BROCK", "GUS", "BYRON", "Equipo 2", "NXT|Arthur", "NXT|Rup", "NXT|amos", "FUT|GeRo", "FUT|MeOw", "FUT|Nowy297", "20250723T173143.000Z")
Set-RowValues $drySeason 25 @("CORDELIUS", "BELLE", "MEEPLE", "JANET", "ALLI", "DOUG", "Equipo 2", "IC|Mebius", "IC|RamaZR", "IC|Nob", "Enraged 💔", "SUP|Tomzy", "SUP|Filippo神", "20250723T173533.000Z")
Set-RowValues $drySeason 26 @("CORDELIUS", "BELLE", "MEEPLE", "JANET", "ALLI", "DOUG", "Equipo 1", "IC|Mebius", "IC|RamaZR", "IC|Nob", "Enraged 💔", "SUP|Tomzy", "SUP|Filippo神", "20250723T173313.000Z")
Set-RowValues $drySeason 27 @("JAE-YONG", "MEEPLE", "BUSTER", "SQUEAK", "GENE", "KAZE", "Equipo 2", "HMB|BosS", "HMB|Lukii", "HMB|Symantec", "TH|iKaoss", "TH|Zhar", "TH|LeNain", "20250723T173941.000Z")
Set-RowValues $drySeason 28 @("JAE-YONG", "MEEPLE", "BUSTER", "SQUEAK", "GENE", "KAZE", "Equipo 2", "HMB|BosS", "HMB|Lukii", "HMB|Symantec", "TH|iKaoss", "TH|Zhar", "TH|LeNain", "20250723T173809.000Z")
Set-RowValues $drySeason 29 @("CARL", "KAZE", "GUS", "HANK", "MEEPLE", "JAE-YONG", "Equipo 2", "HMB|BosS", "HMB|Symantec", "HMB|Lukii", "TH|LeNain", "TH|iKaoss", "TH|Zhar", "20250723T173242.000Z")

# Fix up the "Equipo 1" (winner) cell style on the rows where team 1 won
$drySeason.Range("G4").Copy()
$drySeason.Range("G21").PasteSpecial(-4122)
$drySeason.Range("G23").PasteSpecial(-4122)
$drySeason.Range("G26").PasteSpecial(-4122)

# --- New sheet "Sneaky Fields" (same layout as the other arena sheets) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sneaky = $wb.Worksheets.Add($null, $lastSheet)
$sneaky.Name = "Sneaky Fields"

# Header row (labels + styling) is identical to the other arena sheets
$pitStop.Range("A3:N3").Copy()
$sneaky.Range("A3:N3").PasteSpecial(-4122)
Set-RowValues $sneaky 3 @("B1", "B2", "B3", "B1", "B2", "B3", "Ganador", "Jugador 1", "Jugador 2", "Jugador 3", "Jugador 4", "Jugador 5", "Jugador 6", "Timestamp")

# Data rows 4-7 share the same column styling as a normal "Equipo 2"-win row
$drySeason.Range("A18:N18").Copy()
$sneaky.Range("A4:N7").PasteSpecial(-4122)
Set-RowValues $sneaky 4 @("CORDELIUS", "CHARLIE", "ALLI", "SPIKE", "BONNIE", "BUSTER", "Equipo 2", "IC|Mebius", "IC|RamaZR", "IC|Nob", "Enraged 💔", "SUP|Filippo神", "SUP|Tomzy", "20250723T175519.000Z")
Set-RowValues $sneaky 5 @("CORDELIUS", "CHARLIE", "ALLI", "SPIKE", "BONNIE", "BUSTER", "Equipo 2", "IC|Mebius", "IC|RamaZR", "IC|Nob", "Enraged 💔", "SUP|Filippo神", "SUP|Tomzy", "20250723T175411.000Z")
Set-RowValues $sneaky 6 @("STU", "MEEPLE", "FRANK", "WILLOW", "EL PRIMO", "LUMI", "Equipo 2", "IC|Mebius", "IC|Nob", "IC|RamaZR", "Enraged 💔", "SUP|Filippo神", "SUP|Tomzy", "20250723T174811.000Z")
Set-RowValues $sneaky 7 @("STU", "MEEPLE", "FRANK", "WILLOW", "EL PRIMO", "LUMI", "Equipo 2", "IC|Mebius", "IC|Nob", "IC|RamaZR", "Enraged 💔", "SUP|Filippo神", "SUP|Tomzy", "20250723T174615.000Z")

